{"js": "const body = context.document.body;\nconst searchResults = body.search(\"recently decided to return to work following a break to raise my family, \", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target text not found\");\n}\n\nsearchResults.items[0].insertText(\n  \"decided to return to work following a break to raise my family, \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# Replace = 2 -> wdReplaceAll\n$find.Execute(\n    \"recently decided to return to work following a break to raise my family, \",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"decided to return to work following a break to raise my family, \",\n    2\n)\n"}
